$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d0090608ac3c5d16d7f9fb7cc22c75ae8036352/e2e/507136a8-379e-44e1-9498-351a84a40d30.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8236edcf6dc0fcf848ebdf351fb3e64e7f420b57/e2e/507136a8-379e-44e1-9498-351a84a40d30.md."
# This host rounds the stored OOXML column width to ColumnWidth + 5/6, so back
# off by 5/6 here to land the saved <col width="..."> attribute on exactly 40.
$targetColWidth = 40 - (5/6)

# --- Overview sheet: row 3 is the 507136a8-... record ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-23 20:50:49"

# --- zh-cn sheet: row 3 is the 507136a8-... record ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-23 20:50:44"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColWidth

# --- de-de sheet: row 3 is the 507136a8-... record ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-23 20:50:49"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColWidth
